$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Date" column (B2:B4) for the RAD test cases with new timestamps
$ws.Range("B2").Value = "Tue Jan 28 22:10:56 EST 2025"
$ws.Range("B3").Value = "Tue Jan 28 22:11:10 EST 2025"
$ws.Range("B4").Value = "Tue Jan 28 22:11:24 EST 2025"
